$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in missing "parent company" (column A) values for a few brand rows,
# and correct Hering's parent (acquired by Grupo Soma).
$ws.Range("A16").Value = "Avon"
$ws.Range("A9").Value = "Nubank"
$ws.Range("A8").Value = "Grupo Boticário"
$ws.Range("A57").Value = "Grupo Soma"

# Left-align the brand name in B12.
$ws.Range("B12").HorizontalAlignment = -4131

# Remove the autofilter from the sheet.
$ws.AutoFilterMode = $false

# Move the active selection to B1.
$ws.Range("B1").Select()
